$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix up row 18 (previously the last row) - correct the Run Time value
$ws.Cells.Item(18, 2).Value = 44806.70102041667

# Append the new row 19 with the latest registration/education sprint run
$ws.Cells.Item(19, 1).Value = "2022-09-06"
$ws.Cells.Item(19, 2).Value = 44810.94112326238
$ws.Cells.Item(19, 3).Value = "edu98"
$ws.Cells.Item(19, 4).Value = 60
$ws.Cells.Item(19, 5).Value = 60
$ws.Cells.Item(19, 6).Value = 0
$ws.Cells.Item(19, 7).Value = 1.16

Write-Output "done"
